# ===================================================================
# Update on 20181223.1435 by YKBKyle
# Add new R graphics/stats notes (title, points, legend, lm, abline,
# mtext, text, gl) to the RCommands sheet, rows 142-170.
# ===================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Workbook window position (best effort; may be a no-op) ----
try { $excel.Left = 1780; $excel.Top = 780 } catch {}

# ---- Populate new rows 142-170 with values ----
$ws.Cells.Item(142,1).Value2 = "Package: graphics"
$ws.Cells.Item(142,2).Value2 = "title()"
$ws.Cells.Item(142,3).Value2 = "with(airquality, plot(Wind, Ozone))`ntitle(main = `"Ozone and Wind in New York City`")"
$ws.Cells.Item(142,4).Value2 = "add a title to the plot"
$ws.Cells.Item(144,3).Value2 = "title(`"Ozone and Wind in New York City`")"
$ws.Cells.Item(145,1).Value2 = "Package: graphics"
$ws.Cells.Item(145,2).Value2 = "points()"
$ws.Cells.Item(145,3).Value2 = "with(airquality, plot(Wind, Ozone, main = `"Ozone and Wind in New York City`"))`nwith(subset(airquality, Month == 5), points(Wind, Ozone, col = `"blue`"))"
$ws.Cells.Item(145,4).Value2 = "create the title with one call with plot function;`nget a subset of airquality whose Month is May, and re-add the point (Wind, Ozone) of May in blue"
$ws.Cells.Item(147,1).Value2 = "Package: graphics"
$ws.Cells.Item(147,2).Value2 = "legend()"
$ws.Cells.Item(147,3).Value2 = "with(airquality, plot(Wind, Ozone, main = `"Ozone and Wind in New York City`", type = `"n`"))`nwith(subset(airquality, Month == 5), points(wind, Ozone, col = `"blue`"))`nwith(subset(airquality, Month != 5), points(Wind, Ozone, col = `"red`"))`nlegend(`"topright`", pch = 1, col = c(`"blue`",`"red`"), legend = c(`"May`", `"Other Months`"))"
$ws.Cells.Item(147,4).Value2 = "type = `"n`" means do not plot anything but just setup the graphic device`nso the scales are set but no points are plotted;`nadd blue points and red points to the gragh;`nadd legends to the graph on the top right, shown as rings, with colors and annotation"
$ws.Cells.Item(151,1).Value2 = "Package: stats"
$ws.Cells.Item(151,2).Value2 = "lm()"
$ws.Cells.Item(151,3).Value2 = "model <- lm(Ozone ~ Wind, airquality)"
$ws.Cells.Item(151,4).Value2 = "do linear regression on airquality`$Ozone and airquality`$Wind, just create this model, no plot"
$ws.Cells.Item(152,1).Value2 = "Package: graphics"
$ws.Cells.Item(152,2).Value2 = "abline()"
$ws.Cells.Item(152,3).Value2 = "with(airquality, plot(Wind, Ozone, main = `"Ozone and Wind in New York City`", pch = 20))`nmodel <- lm(Ozone ~ Wind, airquality)`nabline(model, lwd = 2)"
$ws.Cells.Item(152,4).Value2 = "scatter plot points (Wind, Ozone), points showns as solid circles;`ncreate a linear regression on (Wind, Ozone);`nadd the regression line to the plot"
$ws.Cells.Item(154,1).Value2 = "Package: graphics"
$ws.Cells.Item(154,2).Value2 = "par()"
$ws.Cells.Item(154,3).Value2 = "par(`"lty`")"
$ws.Cells.Item(154,4).Value2 = "return the default line type for global graphics parameters"
$ws.Cells.Item(155,3).Value2 = "par(mar = c(2,2,2,2))"
$ws.Cells.Item(155,4).Value2 = "set the margin of 4 sides to be all 2, the default is c(4,4,2,1)"
$ws.Cells.Item(156,3).Value2 = "par(mfrow = c(1, 2))`nwith(airquality, {`nplot(Wind, Ozone, main = `"Ozone and Wind`")`nplot(Solar.R, Ozone, main = `"Ozone and Solar Radiation`")`n})"
$ws.Cells.Item(156,4).Value2 = "create 2 graphic devices, placed in 1 row of 2 columns;`nscatter plot (Wind, Ozone) on the left graphic device;`nscatter plot (Solar.R, Ozone) on the right graphic device"
$ws.Cells.Item(161,1).Value2 = "Package: graphics"
$ws.Cells.Item(161,2).Value2 = "mtext()"
$ws.Cells.Item(161,3).Value2 = "par(mfrow = c(1,3), mar = c(4,4,2,1), oma = c(0,0,2,0))`nwith(airquality, {`nplot(Wind, Ozone, main = `"Ozone and Wind`")`nplot(Solar.R, Ozone, main = `"Ozone and Solar Radiation`")`nplot(Temp, Ozone, main = `"Ozone and Temperature`")`nmtext(`"Ozone and Weather in New York City`", outer = TRUE)`n})"
$ws.Cells.Item(161,4).Value2 = "set mar to adjust the distance between the 3 plots;`nset oma to adjust the outer margin, default: c(0,0,0,0). Oma = c(0,0,2,0) leaves space on the top`nmtext add outer title for all 3 plots"
$ws.Cells.Item(168,1).Value2 = "Package: graphics"
$ws.Cells.Item(168,2).Value2 = "text()"
$ws.Cells.Item(168,3).Value2 = "x <- rnorm(100); y <- rnorm(100); plot(x,y,pch=20); text(0,0, `"label`")"
$ws.Cells.Item(168,4).Value2 = "add text: `"label`" in the location of (0,0) in the plot"
$ws.Cells.Item(169,1).Value2 = "Package: base"
$ws.Cells.Item(169,2).Value2 = "gl()"
$ws.Cells.Item(169,3).Value2 = "gl(2,3)"
$ws.Cells.Item(169,4).Value2 = "create a factor object with 2 levels, specified as 1 and 2, each repeats for 3 times"
$ws.Cells.Item(170,3).Value2 = "gl(2,3, labels = c(`"Male`",`"Female`"))"
$ws.Cells.Item(170,4).Value2 = "create a factor object with 2 levels, specified as Male and Female, each repeats for 3 times"

# ---- Apply wrap text formatting to long note cells ----
$ws.Cells.Item(142,3).WrapText = $true
$ws.Cells.Item(145,3).WrapText = $true
$ws.Cells.Item(145,4).WrapText = $true
$ws.Cells.Item(147,3).WrapText = $true
$ws.Cells.Item(147,4).WrapText = $true
$ws.Cells.Item(148,3).WrapText = $true
$ws.Cells.Item(149,3).WrapText = $true
$ws.Cells.Item(150,3).WrapText = $true
$ws.Cells.Item(152,3).WrapText = $true
$ws.Cells.Item(152,4).WrapText = $true
$ws.Cells.Item(153,3).WrapText = $true
$ws.Cells.Item(156,3).WrapText = $true
$ws.Cells.Item(156,4).WrapText = $true
$ws.Cells.Item(161,3).WrapText = $true
$ws.Cells.Item(161,4).WrapText = $true
$ws.Cells.Item(162,3).WrapText = $true
$ws.Cells.Item(163,3).WrapText = $true
$ws.Cells.Item(164,3).WrapText = $true
$ws.Cells.Item(165,3).WrapText = $true
$ws.Cells.Item(166,3).WrapText = $true
$ws.Cells.Item(167,3).WrapText = $true

# ---- Normalize row heights back to automatic (removes stray custom heights) ----
for ($r = 142; $r -le 170; $r++) {
  $ws.Rows.Item($r).AutoFit()
}

# ---- Apply the few explicit custom row heights present in the target ----
$ws.Rows.Item(147).RowHeight = 18
$ws.Rows.Item(152).RowHeight = 36
$ws.Rows.Item(161).RowHeight = 18

# ---- Merge cells that belong together as single notes ----
$ws.Range("A142:A144").Merge() | Out-Null
$ws.Range("B142:B144").Merge() | Out-Null
$ws.Range("C142:C143").Merge() | Out-Null
$ws.Range("D142:D144").Merge() | Out-Null
$ws.Range("A145:A146").Merge() | Out-Null
$ws.Range("B145:B146").Merge() | Out-Null
$ws.Range("C145:C146").Merge() | Out-Null
$ws.Range("D145:D146").Merge() | Out-Null
$ws.Range("A147:A150").Merge() | Out-Null
$ws.Range("B147:B150").Merge() | Out-Null
$ws.Range("C147:C150").Merge() | Out-Null
$ws.Range("D147:D150").Merge() | Out-Null
$ws.Range("A152:A153").Merge() | Out-Null
$ws.Range("B152:B153").Merge() | Out-Null
$ws.Range("C152:C153").Merge() | Out-Null
$ws.Range("D152:D153").Merge() | Out-Null
$ws.Range("A154:A160").Merge() | Out-Null
$ws.Range("B154:B160").Merge() | Out-Null
$ws.Range("C156:C160").Merge() | Out-Null
$ws.Range("D156:D160").Merge() | Out-Null
$ws.Range("A161:A167").Merge() | Out-Null
$ws.Range("B161:B167").Merge() | Out-Null
$ws.Range("C161:C167").Merge() | Out-Null
$ws.Range("D161:D167").Merge() | Out-Null
$ws.Range("A169:A170").Merge() | Out-Null
$ws.Range("B169:B170").Merge() | Out-Null

# ---- Widen column C to fit the new, longer note text ----
$ws.Columns.Item(3).ColumnWidth = 93.86

# ---- Selection / scroll position (best effort) ----
$ws.Range("C181").Select()
try {
  $win = $excel.Windows.Item(1)
  $win.ScrollRow = 160
  $win.ScrollColumn = 1
} catch {}
